# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The two worker records in the "Estado de Cuenta" table swap positions
# (PAOLA ANDREA DEVIA DEOSSA moves up to row 16, YULEY MARGARITA ALTAHONA
# SANTOYA moves down to row 17), and YULEY's "Salario Basico" (column G)
# is updated to its new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 -> PAOLA ANDREA DEVIA DEOSSA's record
$ws.Range("C16").Value = 1143343119
$ws.Range("D16").Value = "PAOLA ANDREA DEVIA DEOSSA"
$ws.Range("E16").Value = 1802
$ws.Range("F16").Value = 31320
$ws.Range("G16").Value = 783000

# Row 17 -> YULEY MARGARITA ALTAHONA SANTOYA's record (with updated salary)
$ws.Range("C17").Value = 45592178
$ws.Range("D17").Value = "YULEY MARGARITA ALTAHONA SANTOYA"
$ws.Range("E17").Value = 2003
$ws.Range("F17").Value = 44855
$ws.Range("G17").Value = 1320800
